$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(160).Insert()

$ws.Cells.Item(160, 1).Value = 8
$ws.Cells.Item(160, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44964
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = 100112044
$ws.Cells.Item(160, 7).Value = "Perejil"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 2000
$ws.Cells.Item(160, 11).Value = 2800
$ws.Cells.Item(160, 12).Value = 3000
$ws.Cells.Item(160, 13).Value = 2900
$ws.Cells.Item(160, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(160, 16).Value = 1933
$ws.Cells.Item(160, 17).Value = 1.5
$ws.Cells.Item(160, 18).Value = "Hortaliza"
